# Daily attendance processing - 2026-01-10 04:25:08
# Swap the order of names in the "Recorded By" column (G) so that entries
# listing "dnasr281@gmail.com, System" read "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
